$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 13: latest intraday sampling window added to the BBVA table ---
$ws.Range("A13").Value = "Ticker/Index Name: 29/05/2025 15:30-16:30"

$ws.Range("B13").Value = 10.533529
$ws.Range("C13").Value = 15.760683
$ws.Range("D13").Value = 12.350709
$ws.Range("E13").Value = 6.705165
$ws.Range("F13").Value = 4.350376
$ws.Range("G13").Value = 4.463942
$ws.Range("H13").Value = 3.702788
$ws.Range("I13").Value = "N/A"
$ws.Range("J13").Value = "N/A"

# Give the numeric cells the same muted grey Arial look already used
# elsewhere in the sheet (e.g. B9), just at a smaller 9pt size.
$ws.Range("B9").Copy()
foreach ($addr in @("B13", "D13", "F13", "G13", "H13")) {
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats - keeps the values already written above
}
foreach ($addr in @("B13", "D13", "F13", "G13", "H13")) {
    $ws.Range($addr).Font.Size = 9
}

$ws.Range("J13").Select() | Out-Null
